$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on all changed cells so values are stored as literal
# text (matching the original t="inlineStr" cells) rather than being
# auto-coerced into numbers/percentages by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '312.86'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '1.19%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '38.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.46%'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.90%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07928'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.88%'
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.907'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-2.91%'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '8.267'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-0.25%'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9265'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-0.11%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1207'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-8.97%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1916'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-5.74%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09234'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '5.56%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03346'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-2.57%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09637'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.73%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001367'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-1.77%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005923'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-2.19%'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.537'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.45%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.405'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.10%'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.099'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.44%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3449'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.60%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.285'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '5.69%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.98%'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '3.82%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.04379'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.25%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2.32%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004297'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-6.40%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001220'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-9.60%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02115'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-7.39%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05106'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '1.09%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007640'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.12%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.009117'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-8.33%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1360'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '0.58%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002051'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '3.54%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.008603'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-1.19%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006687'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '1.25%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.29%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-3.35%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.08%'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.29%'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.29%'
